# excelhere/test.xlsx -- "saving tree": append three new item rows (43-45)
# below the existing data, same as the webapp's clickEvent()/save() row
# append. Each new row reuses the blank/inline formatting of the row that
# was previously last (row 41/42), so after writing each row's values we
# paste-special just the *formats* from the row above onto it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @{ Row = 43; Name = "qwd";      Unit = "ee";   Price = "qwd" },
    @{ Row = 44; Name = "ee";       Unit = "ef";   Price = "w" },
    @{ Row = 45; Name = "초코렛";   Unit = "봉지"; Price = "3000" }
)

foreach ($item in $newRows) {
    $n = $item.Row

    # Columns A (no.), B (code), F (qty) and G (amount) are left blank for
    # these rows, just like the existing "quick add" rows above them.
    # A leading apostrophe forces every written cell -- blank or not -- to
    # be stored as text, matching the rest of the sheet's inline strings.
    $ws.Range("A$n").Value2 = "'"
    $ws.Range("B$n").Value2 = "'"
    $ws.Range("C$n").Value2 = $item.Name
    $ws.Range("D$n").Value2 = $item.Unit
    $ws.Range("E$n").Value2 = "'" + $item.Price
    $ws.Range("F$n").Value2 = "'"
    $ws.Range("G$n").Value2 = "'"

    # Re-stamp this row's formatting from the row directly above it so it
    # keeps the same (default) style instead of the column's style, and so
    # it doesn't keep the quote-prefix format the apostrophe trick adds.
    $prev = $n - 1
    $ws.Range("A${prev}:G${prev}").Copy() | Out-Null
    $ws.Range("A${n}:G${n}").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false
